$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Test_1" -> "Test_2" (new track entry: FSATA)
$ws.Range("A7").Value = "Test_2"

# Updated raw measurements for the new track / run
$ws.Range("E7").Value = 57.009999999999962
$ws.Range("G7").Value = 5.4099999999999291
$ws.Range("I7").Value = 4.42999999999995

# Match the active cell selection left by the author
$ws.Range("E7").Select()
